$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the existing
# header style (bold, bordered, centered) from H1 so the new headers
# match the look of the other column headers.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data for columns I (I0) and J (IF), keyed by row number.
$data = @{
    2  = @(8, 8)
    3  = @(6, 6)
    4  = @(9, 9)
    5  = @(8, 8)
    6  = @(8, 8)
    7  = @(7, 8)
    8  = @(6, 6)
    9  = @(7, 7)
    10 = @(6, 6)
    11 = @(7, 7)
    12 = @(5, 6)
    13 = @(8, 8)
    14 = @(6, 7)
    15 = @(8, 8)
    16 = @(6, 6)
    17 = @(6, 6)
    18 = @(3, 3)
    19 = @(6, 6)
    20 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
